# Product Sheet update
#
# - Renames the worksheet from "Sheet2" to "Sheet1".
# - Refreshes the five sample product rows (new product codes / brand /
#   sub-category text) and clears the now-unused detail columns I:J.
# - Narrows columns L:N (they were auto "best fit"; now fixed widths).
# - Leaves the active selection on C12 instead of the whole sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

# New product data for rows 2-6 (ProductName / SubCategory / SubChildCategory).
# ProductCode/HSNCode/PackSize/UOM/Quantity/Color/PaintType/FinishType/UnitPrice
# are unchanged, only A, G and H move to the refreshed taxonomy.
$rows = @(
    @{ Row = 2; A = "Test 01"; H = "Tinters" },
    @{ Row = 3; A = "Test 02"; H = "Tinters" },
    @{ Row = 4; A = "Test 03"; H = "Comp-A" },
    @{ Row = 5; A = "Test 04"; H = "Comp-A" },
    @{ Row = 6; A = "Test 05"; H = "Comp-A" }
)

foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("G" + $r.Row).Value = "Jotun"
    $ws.Range("H" + $r.Row).Value = $r.H
}

# SubCategory (I) / SubChildCategory (J) detail cells are no longer populated.
$ws.Range("I2:J6").ClearContents() | Out-Null

# Resize columns L (PackSize detail), M (UOM detail) and N (Quantity detail).
# ColumnWidth is expressed in characters and gets snapped to the workbook's
# pixel grid on save, same as typing a width into Excel's "Column Width"
# dialog - these inputs are the values that land closest to the intended
# 9 / 9.57 / 8.43 character widths.
$ws.Columns.Item(12).ColumnWidth = 8.2
$ws.Columns.Item(13).ColumnWidth = 8.7
$ws.Columns.Item(14).ColumnWidth = 7.7

# Active cell/selection ends up on C12.
$ws.Range("C12").Select() | Out-Null
